{"js": "// Office.js (Word JavaScript API) script\n// Applies the \"Consciousness\" -> \"Arts\" essay rewrite described by the diff:\n//  - Title, author name and e-mail are changed\n//  - Every sentence of the body + summary paragraphs is replaced with new wording\n//  - A couple of new sentences are appended in two places\n//  - A new empty paragraph is added at the very end of the document body\n\nconst body = context.document.body;\n\n// Helper: replace the first (and only) occurrence of `find` with `replacement`,\n// preserving the formatting of the run(s) containing `find` (search + replace\n// keeps the character formatting of the text it overwrites).\nasync function replaceOnce(find, replacement) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + find);\n  }\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1. Title\nawait replaceOnce(\n  \"Unraveling the Enigma of Consciousness: A Journey of Discovery\",\n  \"The Enriching World of Arts: A Journey Through Expression and Creativity\"\n);\n\n// 2. Author name\nawait replaceOnce(\"Vivian Lam\", \"Abigail Kent\");\n\n// 3. E-mail address (local part before the literal \".com\")\nawait replaceOnce(\"vivianlam@protonmail\", \"abigail.kent25@xyz\");\n\n// 4. Opening paragraph (four sentences + two paragraphs separated by line breaks)\nawait replaceOnce(\n  \"Consciousness, the enigmatic phenomenon at the core of our existence, has long captivated philosophers, scientists, and artists alike\",\n  \"As we embark on a voyage through the realm of Arts, we step into a world where imagination reigns supreme, where the boundaries of reality blur, and where emotions find their voice\"\n);\n\nawait replaceOnce(\n  \" As we navigate the complexities of the physical world, our perception of reality and sense of self emerge as remarkable features that define our human experience\",\n  \" Arts, in its myriad forms, serves as a mirror to society, reflecting the human experience in all its complexities and nuances\"\n);\n\nawait replaceOnce(\n  \" Yet, the nature of consciousness remains shrouded in mystery\",\n  \" From the strokes of a paintbrush on a canvas to the flowing melodies of a symphony, art transcends the limitations of language, inviting us to connect with ourselves, with others, and with the world around us\"\n);\n\nawait replaceOnce(\n  \" With each new discovery, we catch tantalizing glimpses into this intricate tapestry of perception, emotion, and thought\",\n  \" In this essay, we will delve into the captivating power of Arts, exploring its multifaceted dimensions and unraveling the profound impact it has on our lives\"\n);\n\n// 5. Second \"paragraph\" (science) inside the same Word paragraph, after the <w:br/>\nawait replaceOnce(\n  \"In the realm of science, researchers are undertaking groundbreaking studies on neural networks and brain activity, seeking to decode the physical mechanisms that underpin consciousness\",\n  \"Within the tapestry of human existence, art weaves its way seamlessly, becoming an integral part of our history, culture, and identity\"\n);\n\nawait replaceOnce(\n  \" Explorations into altered states of consciousness, such as dreams and meditative practices, provide unique insights into the plasticity of our mental landscapes\",\n  \" From the cave paintings of ancient civilizations to the modern masterpieces adorning museums, art serves as a living testament to our collective memory, bridging the gap between generations and providing invaluable insights into the human condition\"\n);\n\nawait replaceOnce(\n  \" The study of consciousness not only enriches our understanding of ourselves but also holds profound implications for our comprehension of artificial intelligence and its potential to emulate human cognition\",\n  \" Whether it's through the lens of literature, the stage of theater, or the notes of music, art immortalizes the triumphs and tribulations of humanity, capturing the essence of what it means to be human\"\n);\n\n// 6. Third \"paragraph\" (art connections), after the second <w:br/>\nawait replaceOnce(\n  \"As we delved into the depths of consciousness, we uncover profound connections to our artistic and cultural expressions\",\n  \"The realm of art transcends the boundaries of mere aesthetics; it possesses the transformative power to shape our perceptions, provoke thought, and inspire action\"\n);\n\nawait replaceOnce(\n  \" From literature's exploration of inner turmoil to music's ability to evoke an emotional response, we find echoes of our conscious experiences reflected in the works of great artists\",\n  \" Art has the ability to challenge societal norms, question established ideologies, and ignite movements for change\"\n);\n\n// This sentence also gains two brand-new trailing sentences.\nawait replaceOnce(\n  \" These explorations transcending disciplinary boundaries offer a multi-faceted perspective on the multifaceted nature of consciousness\",\n  \" It can educate, inform, and empower, giving voice to the marginalized and shedding light on pressing issues.\" +\n    \" By engaging with works of art, we become more empathetic, more tolerant, and more aware of the world around us, fostering a sense of global citizenship and interconnectedness\"\n);\n\n// 7. Summary paragraph\nawait replaceOnce(\n  \"This essay delved into the enigmatic nature of consciousness, weaving together scientific investigations, philosophical contemplations, and artistic representations\",\n  \"In the realm of Arts, we find a world where imagination and creativity flourish, where expression transcends the limitations of language, and where emotions find their voice\"\n);\n\nawait replaceOnce(\n  \" As we continue to unravel the complexities of our conscious experience, we gain a deeper appreciation for the richness and mystery of our own existence\",\n  \" Art serves as a mirror to society, reflecting the human experience in all its complexities and nuances\"\n);\n\n// Final summary sentence is replaced and extended with three brand-new sentences.\nawait replaceOnce(\n  \" The pursuit of understanding consciousness remains an ongoing journey, beckoning us to explore the vast landscapes of our interconnected minds\",\n  \" It has the power to shape our perceptions, provoke thought, and inspire action, becoming an integral part of our history, culture, and identity.\" +\n    \" Art educates, informs, and empowers, fostering empathy, tolerance, and a sense of global citizenship.\" +\n    \" As we engage with works of art, we embark on a transformative journey, enriching our lives and deepening our understanding of ourselves, others, and the world we inhabit\"\n);\n\n// 8. A new, empty paragraph is appended at the very end of the document body.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script\n# Applies the \"Consciousness\" -> \"Arts\" essay rewrite described by the diff:\n#  - Title, author name and e-mail are changed\n#  - Every sentence of the body + summary paragraphs is replaced with new wording\n#  - A couple of new sentences are appended in two places\n#  - A new empty paragraph is added at the very end of the document body\n\n# Keep Word from \"smart quoting\" straight apostrophes in the replacement text\n# (defensive; the Replace-Once helper below sets Range.Text directly, which\n# does not go through AutoFormat/AutoCorrect, but belt-and-braces is cheap).\ntry { $word.Options.AutoFormatAsYouTypeReplaceQuotes = $false } catch {}\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once($findText, $replaceText) {\n    # Locate the single occurrence of $findText and overwrite it by setting\n    # the found Range's .Text directly -- this swaps the run's w:t content\n    # in place (keeping its rPr/formatting) without going through the\n    # \"typed replacement\" AutoCorrect/AutoFormat pipeline (so straight\n    # apostrophes/quotes in $replaceText are not converted to curly ones).\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $found = $find.Execute()\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n    $find.Parent.Text = $replaceText\n}\n\n# 1. Title\nReplace-Once \"Unraveling the Enigma of Consciousness: A Journey of Discovery\" \"The Enriching World of Arts: A Journey Through Expression and Creativity\"\n\n# 2. Author name\nReplace-Once \"Vivian Lam\" \"Abigail Kent\"\n\n# 3. E-mail address (local part before the literal \".com\")\nReplace-Once \"vivianlam@protonmail\" \"abigail.kent25@xyz\"\n\n# 4. Opening paragraph (four sentences + two paragraphs separated by line breaks)\nReplace-Once \"Consciousness, the enigmatic phenomenon at the core of our existence, has long captivated philosophers, scientists, and artists alike\" \"As we embark on a voyage through the realm of Arts, we step into a world where imagination reigns supreme, where the boundaries of reality blur, and where emotions find their voice\"\n\nReplace-Once \" As we navigate the complexities of the physical world, our perception of reality and sense of self emerge as remarkable features that define our human experience\" \" Arts, in its myriad forms, serves as a mirror to society, reflecting the human experience in all its complexities and nuances\"\n\nReplace-Once \" Yet, the nature of consciousness remains shrouded in mystery\" \" From the strokes of a paintbrush on a canvas to the flowing melodies of a symphony, art transcends the limitations of language, inviting us to connect with ourselves, with others, and with the world around us\"\n\nReplace-Once \" With each new discovery, we catch tantalizing glimpses into this intricate tapestry of perception, emotion, and thought\" \" In this essay, we will delve into the captivating power of Arts, exploring its multifaceted dimensions and unraveling the profound impact it has on our lives\"\n\n# 5. Second \"paragraph\" (science) inside the same Word paragraph, after the line break\nReplace-Once \"In the realm of science, researchers are undertaking groundbreaking studies on neural networks and brain activity, seeking to decode the physical mechanisms that underpin consciousness\" \"Within the tapestry of human existence, art weaves its way seamlessly, becoming an integral part of our history, culture, and identity\"\n\nReplace-Once \" Explorations into altered states of consciousness, such as dreams and meditative practices, provide unique insights into the plasticity of our mental landscapes\" \" From the cave paintings of ancient civilizations to the modern masterpieces adorning museums, art serves as a living testament to our collective memory, bridging the gap between generations and providing invaluable insights into the human condition\"\n\nReplace-Once \" The study of consciousness not only enriches our understanding of ourselves but also holds profound implications for our comprehension of artificial intelligence and its potential to emulate human cognition\" \" Whether it's through the lens of literature, the stage of theater, or the notes of music, art immortalizes the triumphs and tribulations of humanity, capturing the essence of what it means to be human\"\n\n# 6. Third \"paragraph\" (art connections), after the second line break\nReplace-Once \"As we delved into the depths of consciousness, we uncover profound connections to our artistic and cultural expressions\" \"The realm of art transcends the boundaries of mere aesthetics; it possesses the transformative power to shape our perceptions, provoke thought, and inspire action\"\n\nReplace-Once \" From literature's exploration of inner turmoil to music's ability to evoke an emotional response, we find echoes of our conscious experiences reflected in the works of great artists\" \" Art has the ability to challenge societal norms, question established ideologies, and ignite movements for change\"\n\n# This sentence also gains two brand-new trailing sentences.\nReplace-Once \" These explorations transcending disciplinary boundaries offer a multi-faceted perspective on the multifaceted nature of consciousness\" \" It can educate, inform, and empower, giving voice to the marginalized and shedding light on pressing issues. By engaging with works of art, we become more empathetic, more tolerant, and more aware of the world around us, fostering a sense of global citizenship and interconnectedness\"\n\n# 7. Summary paragraph\nReplace-Once \"This essay delved into the enigmatic nature of consciousness, weaving together scientific investigations, philosophical contemplations, and artistic representations\" \"In the realm of Arts, we find a world where imagination and creativity flourish, where expression transcends the limitations of language, and where emotions find their voice\"\n\nReplace-Once \" As we continue to unravel the complexities of our conscious experience, we gain a deeper appreciation for the richness and mystery of our own existence\" \" Art serves as a mirror to society, reflecting the human experience in all its complexities and nuances\"\n\n# Final summary sentence is replaced and extended with three brand-new sentences.\nReplace-Once \" The pursuit of understanding consciousness remains an ongoing journey, beckoning us to explore the vast landscapes of our interconnected minds\" \" It has the power to shape our perceptions, provoke thought, and inspire action, becoming an integral part of our history, culture, and identity. Art educates, informs, and empowers, fostering empathy, tolerance, and a sense of global citizenship. As we engage with works of art, we embark on a transformative journey, enriching our lives and deepening our understanding of ourselves, others, and the world we inhabit\"\n\n# 8. A new, empty paragraph is appended at the very end of the document body.\n$endRange = $d.Content\n$endRange.Collapse(0) | Out-Null\n$endRange.InsertParagraphAfter() | Out-Null\n"}
